# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.385.50"
$ws.Range("E2").Value = "  +3.82%  "
$ws.Range("D3").Value = "'3.246.64"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'578.19"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "'182.04"
$ws.Range("E6").Value = "  +6.83%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  -3.63%  "
$ws.Range("D9").Value = "'3.245.94"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("E10").Value = "  +6.02%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("D12").Value = "'0.414"
$ws.Range("E12").Value = "  +5.07%  "
$ws.Range("D13").Value = "'3.812.75"
$ws.Range("E13").Value = "  +2.77%  "
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'28.45"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "'67.367.29"
$ws.Range("E16").Value = "  +3.89%  "
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").Value = "'3.254.92"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").Value = "'376.64"
$ws.Range("E21").Value = "  +5.57%  "
$ws.Range("D22").Value = "'7.61"
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'71.24"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'0.511"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("D27").Value = "'9.59"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "'5.79"
$ws.Range("E30").Value = "  +8.27%  "
$ws.Range("D31").Value = "'1.97"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").Value = "'22.68"
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.77%  "
$ws.Range("D35").Value = "'6.92"
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("D36").Value = "'163.57"
$ws.Range("E36").Value = "  +5.58%  "
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("D38").Value = "'0.850"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").Value = "'1.85"
$ws.Range("E39").Value = "  +4.91%  "
$ws.Range("E40").Value = "  +13.67%  "
$ws.Range("D41").Value = "'26.73"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").Value = "'4.61"
$ws.Range("E42").Value = "  +10.24%  "
$ws.Range("D43").Value = "'2.59"
$ws.Range("E43").Value = "  +4.75%  "
$ws.Range("D44").Value = "'358.36"
$ws.Range("E44").Value = "  +10.47%  "
$ws.Range("D45").Value = "'2.722.54"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("D46").Value = "'25.59"
$ws.Range("E46").Value = "  +6.01%  "
$ws.Range("D47").Value = "'40.82"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").Value = "'0.0679"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").Value = "'0.0279"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("E50").Value = "  +6.50%  "
$ws.Range("E51").Value = "  -0.78%  "
